$d = $word.ActiveDocument

# Word enumeration constants
$wdHeaderFooterPrimary = 1
$wdAlignParagraphCenter = 1

# Add a header (to the document's first/only section) containing the
# questionnaire number, centered, in 12pt Arial - matching the "Header"
# paragraph style already defined in the template.
$section = $d.Sections.First
$header = $section.Headers.Item($wdHeaderFooterPrimary)

$range = $header.Range
$range.InsertAfter("Questionnaire 6")
$range.Style = "Header"
$range.ParagraphFormat.Alignment = $wdAlignParagraphCenter

# Re-select just the inserted text (excluding the trailing paragraph
# mark) so the run-level font formatting doesn't leak onto the pilcrow.
$textRange = $range.Duplicate
$textRange.End = $textRange.End - 1
$textRange.Font.Name = "Arial"
$textRange.Font.Size = 12
